$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1400.7727
$ws.Range("I137").Value = 1251.7693
$ws.Range("J137").Value = 1616
$ws.Range("K137").Value = 3755.3079
$ws.Range("L137").Value = 4848
$ws.Range("M137").Value = -1205.3079
$ws.Range("N137").Value = -9948

$ws.Range("H138").Value = 3118.5398
$ws.Range("I138").Value = 1803.8125
$ws.Range("K138").Value = 5411.4375
$ws.Range("M138").Value = -271.4375

$ws.Range("H140").Value = 37000
$ws.Range("J140").Value = 37000
$ws.Range("L140").Value = 37000
$ws.Range("N140").Value = -47360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 823.7083
$ws.Range("I2").Value = 680.8125
$ws.Range("J2").Value = 1109.5
$ws.Range("K2").Value = 680.8125
$ws.Range("L2").Value = 1109.5
$ws.Range("M2").Value = -567.8125
$ws.Range("N2").Value = -1335.5

$ws.Range("H32").Value = 6086.5
$ws.Range("I32").Value = 4082.1838
$ws.Range("J32").Value = 19500
$ws.Range("K32").Value = 4082.1838
$ws.Range("L32").Value = 19500
$ws.Range("M32").Value = -3795.1838
$ws.Range("N32").Value = -20074

$ws.Range("H61").Value = 2097.9375
$ws.Range("I61").Value = 989
$ws.Range("J61").Value = 3523.7144
$ws.Range("K61").Value = 989
$ws.Range("L61").Value = 3523.7144
$ws.Range("M61").Value = -777
$ws.Range("N61").Value = -3947.7144

$ws.Range("H74").Value = 2571.8809
$ws.Range("I74").Value = 2599.1785
$ws.Range("J74").Value = 2517.2856
$ws.Range("K74").Value = 2599.1785
$ws.Range("L74").Value = 2517.2856
$ws.Range("M74").Value = -1725.1785
$ws.Range("N74").Value = -4265.2856

$ws.Range("H77").Value = 2571.8809
$ws.Range("I77").Value = 2599.1785
$ws.Range("J77").Value = 2517.2856
$ws.Range("K77").Value = 12995.8925
$ws.Range("L77").Value = 12586.428
$ws.Range("M77").Value = -8627.8925
$ws.Range("N77").Value = -21322.428

$ws.Range("H116").Value = 823.7083
$ws.Range("I116").Value = 680.8125
$ws.Range("J116").Value = 1109.5
$ws.Range("K116").Value = 680.8125
$ws.Range("L116").Value = 1109.5
$ws.Range("M116").Value = 1613.1875
$ws.Range("N116").Value = -5697.5

$ws.Range("H136").Value = 2097.9375
$ws.Range("I136").Value = 989
$ws.Range("J136").Value = 3523.7144
$ws.Range("K136").Value = 2967
$ws.Range("L136").Value = 10571.1432
$ws.Range("M136").Value = -417
$ws.Range("N136").Value = -15671.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 823.7083
$ws.Range("I3").Value = 680.8125
$ws.Range("J3").Value = 1109.5
$ws.Range("K3").Value = 680.8125
$ws.Range("L3").Value = 1109.5
$ws.Range("M3").Value = -566.8125
$ws.Range("N3").Value = -1337.5

$ws.Range("H20").Value = 959.4231
$ws.Range("I20").Value = 649.7222
$ws.Range("J20").Value = 1656.25
$ws.Range("K20").Value = 649.7222
$ws.Range("L20").Value = 1656.25
$ws.Range("M20").Value = -402.7222
$ws.Range("N20").Value = -2150.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2502.2666
$ws.Range("I31").Value = 1163.625
$ws.Range("J31").Value = 3240.8276
$ws.Range("K31").Value = 1163.625
$ws.Range("L31").Value = 3240.8276
$ws.Range("M31").Value = -868.625
$ws.Range("N31").Value = -3830.8276

$ws.Range("H34").Value = 2502.2666
$ws.Range("I34").Value = 1163.625
$ws.Range("J34").Value = 3240.8276
$ws.Range("K34").Value = 1163.625
$ws.Range("L34").Value = 3240.8276
$ws.Range("M34").Value = -961.625
$ws.Range("N34").Value = -3644.8276

$ws.Range("H36").Value = 9000
$ws.Range("J36").Value = 9000
$ws.Range("L36").Value = 9000
$ws.Range("N36").Value = -9776

$ws.Range("H40").Value = 9000
$ws.Range("J40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("N40").Value = -9320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 649.6491
$ws.Range("I5").Value = 385.6111
$ws.Range("J5").Value = 1102.2858
$ws.Range("K5").Value = 1156.8333
$ws.Range("L5").Value = 3306.8574
$ws.Range("M5").Value = -1044.8333
$ws.Range("N5").Value = -3530.8574

$ws.Range("H107").Value = 182.71428
$ws.Range("I107").Value = 182.71428
$ws.Range("K107").Value = 548.14284
$ws.Range("M107").Value = 1371.85716

$ws.Range("H113").Value = 1107.3773
$ws.Range("I113").Value = 434
$ws.Range("J113").Value = 1585.258
$ws.Range("K113").Value = 1302
$ws.Range("L113").Value = 4755.774
$ws.Range("M113").Value = 868
$ws.Range("N113").Value = -9095.774000000001

$ws.Range("H131").Value = 803.875
$ws.Range("I131").Value = 389
$ws.Range("J131").Value = 1052.8
$ws.Range("K131").Value = 1167
$ws.Range("L131").Value = 3158.4
$ws.Range("M131").Value = 3873
$ws.Range("N131").Value = -13238.4

$ws.Range("H135").Value = 649.6491
$ws.Range("I135").Value = 385.6111
$ws.Range("J135").Value = 1102.2858
$ws.Range("K135").Value = 3470.4999
$ws.Range("L135").Value = 9920.572200000001
$ws.Range("M135").Value = -935.4999000000003
$ws.Range("N135").Value = -14990.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 2502
$ws.Range("I22").Value = 2733
$ws.Range("J22").Value = 1809
$ws.Range("K22").Value = 2733
$ws.Range("L22").Value = 1809
$ws.Range("M22").Value = -2204
$ws.Range("N22").Value = -2867

$ws.Range("H102").Value = 1726.7693
$ws.Range("I102").Value = 1604.8
$ws.Range("J102").Value = 2133.3333
$ws.Range("K102").Value = 1604.8
$ws.Range("L102").Value = 2133.3333
$ws.Range("M102").Value = 17.20000000000005
$ws.Range("N102").Value = -5377.3333

$ws.Range("H135").Value = 21975
$ws.Range("J135").Value = 21975
$ws.Range("L135").Value = 21975
$ws.Range("N135").Value = -32115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1804.1765
$ws.Range("I93").Value = 1949.875
$ws.Range("J93").Value = 1674.6666
$ws.Range("K93").Value = 1949.875
$ws.Range("L93").Value = 1674.6666
$ws.Range("M93").Value = -701.875
$ws.Range("N93").Value = -4170.6666

$ws.Range("H132").Value = 10645046
$ws.Range("I132").Value = 21749772
$ws.Range("K132").Value = 65249316
$ws.Range("M132").Value = -65246786

$ws.Range("H136").Value = 7808.909
$ws.Range("I136").Value = 24862.4
$ws.Range("K136").Value = 74587.20000000001
$ws.Range("M136").Value = -72037.20000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H81").Value = 40002892
$ws.Range("I81").Value = 200002740
$ws.Range("J81").Value = 2930.4
$ws.Range("K81").Value = 400005480
$ws.Range("L81").Value = 5860.8
$ws.Range("M81").Value = -400004419
$ws.Range("N81").Value = -7982.8

$ws.Range("H84").Value = 40002892
$ws.Range("I84").Value = 200002740
$ws.Range("J84").Value = 2930.4
$ws.Range("K84").Value = 2000027400
$ws.Range("L84").Value = 29304
$ws.Range("M84").Value = -2000022096
$ws.Range("N84").Value = -39912
